# BESensorLessParams - "Works CL at BE"
# Applies the authored edit to the DefaultPars sheet:
#  - Rename/rescale the current-control anti-windup parameter description (row Ind=19)
#  - Update a couple of calibration values (Ind=30, Ind=105)
#  - Insert a new block of 8 "SLPars.Pars6Step.*" parameters (Ind=134..141)
#  - Refresh the selection to cover the new used range

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DefaultPars")

# --- Row 14 (Ind=19, ClaControlPars.KAWUCur): value + description change ---
$ws.Range("E14").Value = 0.95
$ws.Range("F14").Value = "Anti windup scale for current control"

# --- Row 25 (Ind=30, ClaControlPars.Vdc2Bit2Volt): calibration value update ---
$ws.Range("E25").Value = 0.4118

# --- Row 52 (Ind=105, SysState.StepperCurrent.StaticCurrent): range/value update ---
$ws.Range("D52").Value = 18
$ws.Range("E52").Value = 14

# --- Insert 8 new rows for the 6-step sensorless transition parameters ---
$ws.Rows("78:85").Insert()

$ws.Cells.Item(78, 1).Value = 134
$ws.Cells.Item(78, 2).Value = "SLPars.Pars6Step.TransitionTime"
$ws.Cells.Item(78, 3).Value = 0.0000001
$ws.Cells.Item(78, 4).Value = 0.01
$ws.Cells.Item(78, 5).Value = 0.0013
$ws.Cells.Item(78, 6).Value = "Time for transition allowance on 6-step event"

$ws.Cells.Item(79, 1).Value = 135
$ws.Cells.Item(79, 2).Value = "SLPars.Pars6Step.SummingTime"
$ws.Cells.Item(79, 3).Value = 0.0000001
$ws.Cells.Item(79, 4).Value = 0.01
$ws.Cells.Item(79, 5).Value = 0.0013
$ws.Cells.Item(79, 6).Value = "Time for summing voltage / current for R estimate"

$ws.Cells.Item(80, 1).Value = 136
$ws.Cells.Item(80, 2).Value = "SLPars.Pars6Step.MinimumCur4RCalc"
$ws.Cells.Item(80, 3).Value = 0.0000001
$ws.Cells.Item(80, 4).Value = 100
$ws.Cells.Item(80, 5).Value = 5
$ws.Cells.Item(80, 6).Value = "Minimum current step for a well defined R evaluation"

$ws.Cells.Item(81, 1).Value = 137
$ws.Cells.Item(81, 2).Value = "SLPars.Pars6Step.OpenLoopCurDiDtMax"
$ws.Cells.Item(81, 3).Value = 0.0000001
$ws.Cells.Item(81, 4).Value = 1000
$ws.Cells.Item(81, 5).Value = 12
$ws.Cells.Item(81, 6).Value = "Maximum current rise rate for open loop mode"

$ws.Cells.Item(82, 1).Value = 138
$ws.Cells.Item(82, 2).Value = "SLPars.Pars6Step.MaxStepTime"
$ws.Cells.Item(82, 3).Value = 0.0000001
$ws.Cells.Item(82, 4).Value = 1000
$ws.Cells.Item(82, 5).Value = 0.1
$ws.Cells.Item(82, 6).Value = "Maximum step time in 6 step mode"

$ws.Cells.Item(83, 1).Value = 139
$ws.Cells.Item(83, 2).Value = "SLPars.Pars6Step.JOverKT"
$ws.Cells.Item(83, 3).Value = 0.0000001
$ws.Cells.Item(83, 4).Value = 1000
$ws.Cells.Item(83, 5).Value = 0.001
$ws.Cells.Item(83, 6).Value = "Plant dynamics normalizer"

$ws.Cells.Item(84, 1).Value = 140
$ws.Cells.Item(84, 2).Value = "SLPars.Pars6Step.Har3Phase"
$ws.Cells.Item(84, 3).Value = -6.3
$ws.Cells.Item(84, 4).Value = 6.3
$ws.Cells.Item(84, 5).Value = 0
$ws.Cells.Item(84, 6).Value = "Third harmonic angle correction phase"

$ws.Cells.Item(85, 1).Value = 141
$ws.Cells.Item(85, 2).Value = "SLPars.Pars6Step.Har3Amp"
$ws.Cells.Item(85, 3).Value = 0
$ws.Cells.Item(85, 4).Value = 1
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(85, 6).Value = "Third harmonic angle correction amplitude"

# --- Refresh selection to the new used range ---
$ws.Range("A1:F93").Select()
